$p = $ppt.ActivePresentation

# --- Slide 4: "Name components" ---------------------------------------
# Update the sensor / name-component explanation paragraph (lvl=1 bullet).
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(2).Runs(1).Text = "In the example, “xfmr-6.dmd.inst” is referring to a specific sensor in panel xfmr-6 that publishes “electricity demand” data. The data type is reflected in the later name components (in this case, “electricity – aggregation average”)."

# --- Slide 5: "Hierarchical storage" -----------------------------------
$s5 = $p.Slides.Item(5)
$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2: "Leaf nodes publish aggregated (min, sum, avg, etc) ..."
# Split "aggregated " (within the first run) into "aggregation " while
# also splitting off the leading "Leaf nodes publish " text into its own
# run, matching the three-run structure in the target deck.
$para2 = $tr5.Paragraphs(2)
$offset2 = $para2.Text.IndexOf("aggregated ")
$tr5.Characters($para2.Start + $offset2, 11).Text = "aggregation "

# Paragraph 3: "Non-leaf nodes fetch the aggregated data from ..."
$para3 = $tr5.Paragraphs(3)
$offset3 = $para3.Text.IndexOf("aggregated ")
$tr5.Characters($para3.Start + $offset3, 11).Text = "aggregation "
